$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values per the coinranking.com refresh snapshot.
# Cells are stored as text, so force a text number format before assigning
# the string value (keeps exact formatting, e.g. trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.86"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.413"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05890"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.551"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8102"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9341"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1418"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07427"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03382"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03049"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09344"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.857"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001582"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04667"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005897"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004903"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.563"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.145"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3232"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002296"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03976"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1072"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003001"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009706"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005185"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6702"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002386"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
